# Weekly update: a new week's price record is inserted at the top of the
# "Ajo" (garlic) price block for "Terminal La Palmera de La Serena".
# This pushes the existing records (rows 167-237) down by one row, and the
# last existing record ends up as the new last row (238).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 167, shifting 167:237 -> 168:238
$ws.Rows(167).Insert()

# Populate the newly inserted row 167 with this week's data
$ws.Cells.Item(167, 1).Value  = 8
$ws.Cells.Item(167, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(167, 3).Value  = "Coquimbo"
$ws.Cells.Item(167, 4).Value  = 44636
$ws.Cells.Item(167, 5).Value  = 4
$ws.Cells.Item(167, 6).Value  = 100112003
$ws.Cells.Item(167, 7).Value  = "Ajo"
$ws.Cells.Item(167, 8).Value  = "Chino"
$ws.Cells.Item(167, 9).Value  = "Primera"
$ws.Cells.Item(167, 10).Value = 540
$ws.Cells.Item(167, 11).Value = 19000
$ws.Cells.Item(167, 12).Value = 20000
$ws.Cells.Item(167, 13).Value = 19500
$ws.Cells.Item(167, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(167, 15).Value = "China"
$ws.Cells.Item(167, 16).Value = 1950
$ws.Cells.Item(167, 17).Value = 10
$ws.Cells.Item(167, 18).Value = "Hortaliza"
